$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 "Marking": update Right marking from 4 to 5, Wrong marking from -1 to -1.2
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = -1.2

# Row 12 "Total": update Right total from 96 to 120, Wrong total from -1 to -1.2
$ws.Range("B12").Value = 120
$ws.Range("C12").Value = -1.2

# Update the Max column's text representation
$ws.Range("E12").Value = "118.8/140"
